$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Bdnf"
$ws.Cells.Item(2, 3).Value = "Sort1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.03254066666666667
$ws.Cells.Item(2, 8).Value = 0.097622
$ws.Cells.Item(2, 9).Value = 0.02633076257175775
$ws.Cells.Item(2, 10).Value = 0.03898291886711458
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.9703876666666668
$ws.Cells.Item(2, 14).Value = 2.911163
$ws.Cells.Item(2, 15).Value = 0.03945299285965207
$ws.Cells.Item(2, 16).Value = 0.04754668824173519
$ws.Cells.Item(2, 17).Value = 0.03157706159844445
$ws.Cells.Item(2, 18).Value = 0.284193554386
$ws.Cells.Item(2, 19).Value = 0.001038827387732752
$ws.Cells.Item(2, 20).Value = 0.001853508690127554

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Bdnf"
$ws.Cells.Item(3, 3).Value = "Sort1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.03254066666666667
$ws.Cells.Item(3, 8).Value = 0.097622
$ws.Cells.Item(3, 9).Value = 0.02633076257175775
$ws.Cells.Item(3, 10).Value = 0.03898291886711458
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.041192666666667
$ws.Cells.Item(3, 14).Value = 3.123578
$ws.Cells.Item(3, 15).Value = 0.04233170747586662
$ws.Cells.Item(3, 16).Value = 0.05101596487889641
$ws.Cells.Item(3, 17).Value = 0.03388110350177778
$ws.Cells.Item(3, 18).Value = 0.304929931516
$ws.Cells.Item(3, 19).Value = 0.001114626138804146
$ws.Cells.Item(3, 20).Value = 0.001988751219801586

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Bdnf"
$ws.Cells.Item(4, 3).Value = "Sort1"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.03254066666666667
$ws.Cells.Item(4, 8).Value = 0.097622
$ws.Cells.Item(4, 9).Value = 0.02633076257175775
$ws.Cells.Item(4, 10).Value = 0.03898291886711458
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 5.170211333333333
$ws.Cells.Item(4, 14).Value = 15.510634
$ws.Cells.Item(4, 15).Value = 0.2102049704707969
$ws.Cells.Item(4, 16).Value = 0.2533280614069559
$ws.Cells.Item(4, 17).Value = 0.1682421235942222
$ws.Cells.Item(4, 18).Value = 1.514179112348
$ws.Cells.Item(4, 19).Value = 0.005534857168869902
$ws.Cells.Item(4, 20).Value = 0.009875467264590782

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Bdnf"
$ws.Cells.Item(5, 3).Value = "Sort1"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.03254066666666667
$ws.Cells.Item(5, 8).Value = 0.097622
$ws.Cells.Item(5, 9).Value = 0.02633076257175775
$ws.Cells.Item(5, 10).Value = 0.03898291886711458
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 4.853575333333334
$ws.Cells.Item(5, 14).Value = 14.560726
$ws.Cells.Item(5, 15).Value = 0.1973315197085667
$ws.Cells.Item(5, 16).Value = 0.2378136503161547
$ws.Cells.Item(5, 17).Value = 0.1579385770635556
$ws.Cells.Item(5, 18).Value = 1.421447193572
$ws.Cells.Item(5, 19).Value = 0.005195889393370405
$ws.Cells.Item(5, 20).Value = 0.009270670235767017

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Bdnf"
$ws.Cells.Item(6, 3).Value = "Sort1"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.03254066666666667
$ws.Cells.Item(6, 8).Value = 0.097622
$ws.Cells.Item(6, 9).Value = 0.02633076257175775
$ws.Cells.Item(6, 10).Value = 0.03898291886711458
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 12.56068
$ws.Cells.Item(6, 14).Value = 25.12136
$ws.Cells.Item(6, 15).Value = 0.5106788094851177
$ws.Cells.Item(6, 16).Value = 0.4102956351562577
$ws.Cells.Item(6, 17).Value = 0.4087329009866667
$ws.Cells.Item(6, 18).Value = 2.45239740592
$ws.Cells.Item(6, 19).Value = 0.01344656248298054
$ws.Cells.Item(6, 20).Value = 0.01599452145682764

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Bdnf"
$ws.Cells.Item(7, 3).Value = "Sort1"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.2033015
$ws.Cells.Item(7, 8).Value = 2.406603
$ws.Cells.Item(7, 9).Value = 0.9736692374282422
$ws.Cells.Item(7, 10).Value = 0.9610170811328854
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.9703876666666668
$ws.Cells.Item(7, 14).Value = 2.911163
$ws.Cells.Item(7, 15).Value = 0.03945299285965207
$ws.Cells.Item(7, 16).Value = 0.04754668824173519
$ws.Cells.Item(7, 17).Value = 1.1676689348815
$ws.Cells.Item(7, 18).Value = 7.006013609289001
$ws.Cells.Item(7, 19).Value = 0.03841416547191932
$ws.Cells.Item(7, 20).Value = 0.04569317955160764

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Bdnf"
$ws.Cells.Item(8, 3).Value = "Sort1"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.2033015
$ws.Cells.Item(8, 8).Value = 2.406603
$ws.Cells.Item(8, 9).Value = 0.9736692374282422
$ws.Cells.Item(8, 10).Value = 0.9610170811328854
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.041192666666667
$ws.Cells.Item(8, 14).Value = 3.123578
$ws.Cells.Item(8, 15).Value = 0.04233170747586662
$ws.Cells.Item(8, 16).Value = 0.05101596487889641
$ws.Cells.Item(8, 17).Value = 1.252868697589
$ws.Cells.Item(8, 18).Value = 7.517212185534
$ws.Cells.Item(8, 19).Value = 0.04121708133706247
$ws.Cells.Item(8, 20).Value = 0.04902721365909483

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Bdnf"
$ws.Cells.Item(9, 3).Value = "Sort1"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.2033015
$ws.Cells.Item(9, 8).Value = 2.406603
$ws.Cells.Item(9, 9).Value = 0.9736692374282422
$ws.Cells.Item(9, 10).Value = 0.9610170811328854
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 5.170211333333333
$ws.Cells.Item(9, 14).Value = 15.510634
$ws.Cells.Item(9, 15).Value = 0.2102049704707969
$ws.Cells.Item(9, 16).Value = 0.2533280614069559
$ws.Cells.Item(9, 17).Value = 6.221323052717
$ws.Cells.Item(9, 18).Value = 37.327938316302
$ws.Cells.Item(9, 19).Value = 0.204670113301927
$ws.Cells.Item(9, 20).Value = 0.2434525941423651

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Bdnf"
$ws.Cells.Item(10, 3).Value = "Sort1"
$ws.Cells.Item(10, 4).Value = "M2"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.2033015
$ws.Cells.Item(10, 8).Value = 2.406603
$ws.Cells.Item(10, 9).Value = 0.9736692374282422
$ws.Cells.Item(10, 10).Value = 0.9610170811328854
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 4.853575333333334
$ws.Cells.Item(10, 14).Value = 14.560726
$ws.Cells.Item(10, 15).Value = 0.1973315197085667
$ws.Cells.Item(10, 16).Value = 0.2378136503161547
$ws.Cells.Item(10, 17).Value = 5.840314478963
$ws.Cells.Item(10, 18).Value = 35.041886873778
$ws.Cells.Item(10, 19).Value = 0.1921356303151963
$ws.Cells.Item(10, 20).Value = 0.2285429800803877

# Row 11
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Bdnf"
$ws.Cells.Item(11, 3).Value = "Sort1"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.2033015
$ws.Cells.Item(11, 8).Value = 2.406603
$ws.Cells.Item(11, 9).Value = 0.9736692374282422
$ws.Cells.Item(11, 10).Value = 0.9610170811328854
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 12.56068
$ws.Cells.Item(11, 14).Value = 25.12136
$ws.Cells.Item(11, 15).Value = 0.5106788094851177
$ws.Cells.Item(11, 16).Value = 0.4102956351562577
$ws.Cells.Item(11, 17).Value = 15.11428508502
$ws.Cells.Item(11, 18).Value = 60.45714034008
$ws.Cells.Item(11, 19).Value = 0.4972322470021371
$ws.Cells.Item(11, 20).Value = 0.3943011136994301
